$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (shifts existing data rows 3-9 down to 4-10).
$ws.Rows(3).Insert()

# Copy the formatting of the row above into the freshly inserted row so the
# new row matches the table's existing bordered/left-aligned style instead of
# picking up the bare column default style.
$ws.Range("A2:G2").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New pin entry: PB6 / DO / ENRF.
$ws.Range("B3").Value = "PB6"
$ws.Range("C3").Value = "DO"
$ws.Range("D3").Value = "ENRF"

# Number the "序号" (index) column for every data row now that there are nine
# of them (1..9).
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9

# Match the saved cursor/selection position recorded in the workbook.
[void]$ws.Range("D14").Select()
